# Update patient identification data on the "Hoja de Ingreso y Egreso" sheet.
# (commit: "para agregar codigo de barras")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient name / file number (row 6) ---
$ws.Range("A6").Value = "PÉREZ"
$ws.Range("C6").Value = "CONTRERAS"
$ws.Range("E6").Value = "WENDY"
$ws.Range("G6").Value = "MARÍA LORENA"
$ws.Range("I6").Value = "/201773486"

# --- Current address (row 8) ---
$ws.Range("A8").Value = "14 CALLE 37-12 ZONA 7"
$ws.Range("D8").Value = "EL RODEO"
$ws.Range("H8").Value = "GUATEMALA"
$ws.Range("J8").Value = "5577 4236"

# --- Date of birth / age / place of birth (row 12) ---
# These two values look like a date / a plain number, so a direct .Value
# assignment would make the engine coerce them into a real date serial or
# a numeric value (and attach a "Text" quote-prefixed number format to the
# cell). Force them in as plain text, then restore the original General
# formatting by pasting the format from an untouched neighboring cell that
# shares the same style.
$ws.Range("A12").Value = "'1986-01-03"
$ws.Range("B12").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F12").Value = "'31"
$ws.Range("G12").Copy()
$ws.Range("F12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H12").Value = "GUATEMALA, GUATEMALA"

# --- No. de cedula / DPI (row 14) ---
$ws.Range("H14").Value = "DPI 2587058760101"

# --- Spouse / address note (row 16) ---
$ws.Range("A16").Value = "NOÉ DAVID NAJERA"
$ws.Range("F16").Value = "IDEM"

# --- Father / Mother names (row 18) ---
$ws.Range("A18").Value = "JORGE PÉREZ JUÁREZ"
$ws.Range("F18").Value = "SILVIA VIDALIA CONTRERAS"

# --- Emergency contact (row 20) ---
$ws.Range("A20").Value = "NOE DAVID NAJERA"
$ws.Range("F20").Value = "ESPOSO"
$ws.Range("H20").Value = "IDEM"
$ws.Range("J20").Value = "5577 5327"

# --- Admission date / time / service (row 24) ---
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").Value = "15:36:6"
$ws.Range("D24").Value = "E.G.O."

$excel.ActiveWorkbook.Save()
